# Refresh the cryptos price ticker (GitHub Actions scheduled update).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Column D values that are
# plain decimals ("1.00", "8.25", ...) need NumberFormat "@" forced before
# the write, otherwise Excel's smart entry would coerce them to numbers and
# drop the trailing zero / swallow the padded look of the original text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.478.60'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.451.48'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.38'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.08'
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("D7").Value = '3.453.06'
$ws.Range("E7").Value = '  +0.85%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.479'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.63'
$ws.Range("E10").Value = '  -0.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.127'
$ws.Range("E11").Value = '  +3.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '4.051.35'
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.76'
$ws.Range("E14").Value = '  +7.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.122'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000174'
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.474.53'
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '61.691.09'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.38'
$ws.Range("E19").Value = '  +7.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.35'
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.43'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '401.49'
$ws.Range("E22").Value = '  +5.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.566'
$ws.Range("E23").Value = '  +2.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.91'
$ws.Range("E24").Value = '  +3.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.997'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  +1.35%  '
$ws.Range("D27").Value = '3.595.49'
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.181'
$ws.Range("E28").Value = '  +2.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.64'
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.50'
$ws.Range("E31").Value = '  -5.93%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.25'
$ws.Range("E32").Value = '  +1.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.19'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '24.02'
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("E36").Value = '  +2.11%  '
$ws.Range("D37").Value = '3.485.51'
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.18'
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.56'
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '167.03'
$ws.Range("E40").Value = '  +0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0795'
$ws.Range("E41").Value = '  +2.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.44'
$ws.Range("E42").Value = '  +2.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.806'
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.53'
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '42.45'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("D48").Value = '2.614.43'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.98'
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.13'
$ws.Range("E51").Value = '  -3.06%  '
